$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates written as plain text (e.g. "04/01/2021"), matching
# the rest of the sheet. Force text formatting first so Excel's COM layer
# doesn't auto-convert the string into a date serial number, then reset the
# cell style back to Normal so no stray number-format style lingers on the
# new cells (keeps them identical in appearance to the existing rows).
$ws.Range("A386:A389").NumberFormat = "@"

$ws.Range("A386").Value = "04/01/2021"
$ws.Range("B386").Value = 0.43
$ws.Range("C386").Value = 0.43
$ws.Range("D386").Value = "Quinta-Feira"

$ws.Range("A387").Value = "04/02/2021"
$ws.Range("B387").Value = 0.47
$ws.Range("C387").Value = 0.47
$ws.Range("D387").Value = "Sexta-feira"

$ws.Range("A388").Value = "04/03/2021"
$ws.Range("B388").Value = 0.47
$ws.Range("C388").Value = 0.46
$ws.Range("D388").Value = "Sábado"

$ws.Range("A389").Value = "04/04/2021"
$ws.Range("B389").Value = 0.51
$ws.Range("C389").Value = 0.5
$ws.Range("D389").Value = "Domingo"

$ws.Range("A386:A389").Style = "Normal"
